# Add newer project-plan entries (rows 35-45) to Sheet2, matching the
# "Add files via upload" commit that appended more tracked work items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 35 - fill in date/type/progress note for an already-present blank row
$ws.Range("A35").Value = [DateTime]"2025-04-02"
$ws.Range("C35").Value = "Not much progress"
$ws.Range("E35").Value = 0

# Row 36
$ws.Range("A36").Value = [DateTime]"2025-04-03"
$ws.Range("C36").Value = "Naming convention corrections  again "
$ws.Range("D36").Value = "reaplace XX with CA"
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = "Corrected again"

# Row 37
$ws.Range("A37").Value = [DateTime]"2025-04-11"
$ws.Range("C37").Value = "Odata customer master createion"
$ws.Range("D37").Value = "Analysis (4) to find 3 BAPI's + Classes and to have a success test data, "
$ws.Range("E37").Value = 15
$ws.Range("E37").Borders.LineStyle = 1
$ws.Range("F37").Value = "It's a big master involoving 3 BAPI and a class API which has deep strucutre"

# Row 38
$ws.Range("A38").Value = [DateTime]"2025-04-12"
$ws.Range("D38").Value = "6 hrs for dummy report cretion ( can't have BAPI test in SE37 and it involved deep strucutre"
$ws.Range("F38").Value = "Took time to understand BP vs KUNNR, and I captured reasonable efforts that any developer can take"

# Row 39
$ws.Range("A39").Value = [DateTime]"2025-04-13"
$ws.Range("D39").Value = "5 hrs for odata coding and testing"

# Row 40
$ws.Range("A40").Value = [DateTime]"2025-04-19"
$ws.Range("C40").Value = "Digitax changes"
$ws.Range("E40").Value = 0.5
$ws.Range("F40").Value = "Most of the work Akash did"

# Row 41
$ws.Range("A41").Value = [DateTime]"2025-04-20"
$ws.Range("C41").Value = "Digitax changes"
$ws.Range("E41").Value = 0.5

# Row 42
$ws.Range("A42").Value = [DateTime]"2025-04-21"
$ws.Range("C42").Value = "Ad-hoc meeting - Late night"
$ws.Range("D42").Value = "Meeting during late night"
$ws.Range("E42").Value = 2
$ws.Range("F42").Value = "Was part of meeting for initial support"

# Row 43
$ws.Range("A43").Value = [DateTime]"2025-04-22"
$ws.Range("C43").Value = "Ad-hoc meeting"
$ws.Range("D43").Value = "Meeting during late evening"
$ws.Range("E43").Value = 1

# Row 44
$ws.Range("A44").Value = [DateTime]"2025-04-23"
$ws.Range("C44").Value = "Ad-hoc meeting"
$ws.Range("D44").Value = "Meeting during office hrs"
$ws.Range("E44").Value = 3

# Row 45
$ws.Range("A45").Value = [DateTime]"2025-04-24"
$ws.Range("C45").Value = "Ad-hoc meeting"
$ws.Range("D45").Value = "Meeting during office hrs"
$ws.Range("E45").Value = 2

$ws.Range("E45").Select()
